$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mirror the "Misc." column (B) values into the "Expenses" column (A)
# for the rows where A was previously blank.
$ws.Range("A6").Value = 200
$ws.Range("A7").Value = 165
$ws.Range("A11").Value = 1000
$ws.Range("A14").Value = 40
$ws.Range("A18").Value = 700

# Update the view: scroll back to top-left, set zoom to 100%, and
# change the active selection to A21:A23.
$ws.Activate()
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A21:A23").Select() | Out-Null
